# 김다희 Enum에 PickaxeType 추가
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New enum table: PickaxeType (rows 109-113), mirrors the layout/format
# of the other small enum tables already on the sheet (e.g. rows 9-23 / 91-105)
$ws.Range("A91:C91").Copy() | Out-Null
$ws.Range("A109:C109").PasteSpecial(-4122) | Out-Null

$ws.Range("A10:C10").Copy() | Out-Null
$ws.Range("A110:C112").PasteSpecial(-4122) | Out-Null

$ws.Range("A23:C23").Copy() | Out-Null
$ws.Range("A113:C113").PasteSpecial(-4122) | Out-Null

$ws.Range("A109").Value = "EnumName"
$ws.Range("B109").Value = "Typename"
$ws.Range("C109").Value = "TypeValue"

$ws.Range("A110").Value = "PickaxeType"
$ws.Range("B110").Value = "Wood"
$ws.Range("C110").Value = 1

$ws.Range("B111").Value = "Blue"
$ws.Range("C111").Value = 2

$ws.Range("B112").Value = "Red"
$ws.Range("C112").Value = 3

$ws.Range("B113").Value = "Black"
$ws.Range("C113").Value = 4

# --- Log row in the change-history mini table (I/J/K/L, row 31)
$ws.Range("I30:L30").Copy() | Out-Null
$ws.Range("I31:L31").PasteSpecial(-4122) | Out-Null

$ws.Range("I31").Value = "09.04 21:50"
$ws.Range("J31").Value = "김다희"
$ws.Range("K31").Value = "PickaxeType"
$ws.Range("L31").Value = "추가"

$excel.CutCopyMode = 0
